$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for specific rows, per repulled source data.
$ws.Range("F4").Value = -10
$ws.Range("F6").Value = -2
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -9
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = -10
$ws.Range("F20").Value = -3
$ws.Range("F23").Value = -9
$ws.Range("F26").Value = -6
$ws.Range("F27").Value = -10
$ws.Range("F30").Value = -8
$ws.Range("F31").Value = -10
$ws.Range("F33").Value = -4
$ws.Range("F37").Value = -1
